# Refreshes the cryptos worksheet price/volume columns (and the two
# swapped Aave / FraxShare rows) to match the scraped GitHub Actions run:
# "Updated cryptos list on Fri Nov 17 13:35:33 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($address, $text) {
    # Leading apostrophe forces Excel to store the value as literal TEXT
    # instead of auto-coercing number-looking strings (e.g. "244.47",
    # "1.967.94", "0.0850") into numeric values - matching the original
    # workbook, where every Coin/Link/Price/Volume cell is a string.
    $ws.Range($address).Value = "'" + $text
}

Set-TextCell 'D2' '36.522.91'
Set-TextCell 'E2' '  -1.11%  '
Set-TextCell 'D3' '1.964.82'
Set-TextCell 'E3' '  -3.41%  '
Set-TextCell 'E4' '  +0.05%  '
Set-TextCell 'D5' '244.47'
Set-TextCell 'E5' '  -1.63%  '
Set-TextCell 'D6' '0.619'
Set-TextCell 'E6' '  -2.80%  '
Set-TextCell 'D7' '58.92'
Set-TextCell 'E7' '  -6.56%  '
Set-TextCell 'E8' '  -0.04%  '
Set-TextCell 'D9' '0.374'
Set-TextCell 'E9' '  -2.89%  '
Set-TextCell 'D10' '55.79'
Set-TextCell 'E10' '  -3.90%  '
Set-TextCell 'D11' '0.0851'
Set-TextCell 'E11' '  +6.68%  '
Set-TextCell 'D12' '0.103'
Set-TextCell 'E12' '  -0.42%  '
Set-TextCell 'D13' '22.15'
Set-TextCell 'E13' '  -3.15%  '
Set-TextCell 'D14' '0.840'
Set-TextCell 'E14' '  -6.42%  '
Set-TextCell 'D15' '2.256.75'
Set-TextCell 'E15' '  -3.33%  '
Set-TextCell 'D16' '13.57'
Set-TextCell 'E16' '  -5.40%  '
Set-TextCell 'D17' '5.35'
Set-TextCell 'E17' '  -3.19%  '
Set-TextCell 'D18' '1.972.78'
Set-TextCell 'E18' '  -3.13%  '
Set-TextCell 'D19' '36.450.63'
Set-TextCell 'E19' '  -1.40%  '
Set-TextCell 'D20' '0.0₃0887'
Set-TextCell 'E20' '  +0.85%  '
Set-TextCell 'D21' '70.45'
Set-TextCell 'E21' '  -2.24%  '
Set-TextCell 'D22' '231.36'
Set-TextCell 'E22' '  -2.01%  '
Set-TextCell 'D23' '5.10'
Set-TextCell 'E23' '  -5.18%  '
Set-TextCell 'D24' '0.999'
Set-TextCell 'E24' '  -0.20%  '
Set-TextCell 'D25' '2.52'
Set-TextCell 'E25' '  -0.33%  '
Set-TextCell 'E26' '  -2.77%  '
Set-TextCell 'D27' '9.57'
Set-TextCell 'E27' '  -1.76%  '
Set-TextCell 'D28' '164.94'
Set-TextCell 'E28' '  +3.50%  '
Set-TextCell 'D29' '19.72'
Set-TextCell 'E29' '  -2.06%  '
Set-TextCell 'E30' '  -8.11%  '
Set-TextCell 'E31' '  -1.90%  '
Set-TextCell 'D32' '1.17'
Set-TextCell 'E32' '  -0.23%  '
Set-TextCell 'D33' '4.77'
Set-TextCell 'E33' '  -5.54%  '
Set-TextCell 'D34' '0.0639'
Set-TextCell 'E34' '  +3.41%  '
Set-TextCell 'D35' '4.37'
Set-TextCell 'E35' '  -3.09%  '
Set-TextCell 'D36' '6.17'
Set-TextCell 'E36' '  -2.13%  '
Set-TextCell 'E37' '  +0.15%  '
Set-TextCell 'E38' '  -1.79%  '
Set-TextCell 'D39' '2.18'
Set-TextCell 'E39' '  -8.55%  '
Set-TextCell 'D40' '2.91'
Set-TextCell 'E40' '  -5.59%  '
Set-TextCell 'E41' '  -0.90%  '
Set-TextCell 'D42' '1.19'
Set-TextCell 'E42' '  -3.82%  '
Set-TextCell 'D43' '2.87'
Set-TextCell 'E43' '  -3.36%  '
Set-TextCell 'E44' '  -1.67%  '
Set-TextCell 'D45' '15.88'
Set-TextCell 'E45' '  -7.21%  '
Set-TextCell 'D46' '1.05'
Set-TextCell 'E46' '  -7.03%  '
Set-TextCell 'B47' 'FraxShare'
Set-TextCell 'C47' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D47' '7.43'
Set-TextCell 'E47' '  -3.31%  '
Set-TextCell 'B48' 'Aave'
Set-TextCell 'C48' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D48' '89.30'
Set-TextCell 'E48' '  -4.39%  '
Set-TextCell 'D49' '1.348.68'
Set-TextCell 'E49' '  -1.32%  '
Set-TextCell 'E50' '  -3.14%  '
Set-TextCell 'D51' '45.10'
Set-TextCell 'E51' '  -0.69%  '
